# Arreglo de pago de siniestros
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHost = "preproducciongestion.segurossura.com.ar"
$newUrl  = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"

# --- Remove old hyperlink(s) so we can recreate them in the right order ---
$ws.Hyperlinks.Delete()

# --- Row 2 ---
$ws.Range("B2").Value = $newUrl
$ws.Range("A2").Value = $newHost
$ws.Range("E2").Value = "Incendio parcial"
$ws.Range("F2").Value = "'0420172007039"
$ws.Range("G2").Value = "Cheque"

# --- Row 3 ---
$ws.Range("B3").Value = $newUrl
$ws.Range("A3").Value = $newHost
$ws.Range("E3").Value = "Parcial"
$ws.Range("F3").Value = "'0420172006736"
$ws.Range("G3").Value = "Transferencia electrónica de fondos"

# --- Row 4 ---
$ws.Range("E4").ClearContents() | Out-Null
$ws.Range("F4").Value = "'0420172007039"
$ws.Range("G4").Value = "Cheque"

# --- Row 12 ---
$ws.Range("B12").Value = " "

# --- Recreate hyperlinks: B2 first (rId1), then B3 (rId2) ---
$ws.Hyperlinks.Add($ws.Range("B2"), $newUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $newUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newUrl) | Out-Null

# Ensure hyperlink cells use the built-in "Hipervínculo" style (as before)
$ws.Range("B2").Style = "Hipervínculo"
$ws.Range("B3").Style = "Hipervínculo"

# --- View state: select B3, no frozen/scrolled top-left cell ---
$ws.Activate()
$ws.Range("B3").Select() | Out-Null
